# Auto-generated Excel COM-interop script to apply scheduled price/profit updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3390.85
$ws.Range("I6").Value = 3237.7222
$ws.Range("K6").Value = 9713.1666
$ws.Range("M6").Value = -9601.1666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4227.75
$ws.Range("I62").Value = 3970.3333
$ws.Range("K62").Value = 3970.3333
$ws.Range("M62").Value = -3346.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 4227.75
$ws.Range("I65").Value = 3970.3333
$ws.Range("K65").Value = 19851.6665
$ws.Range("M65").Value = -16731.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 884
$ws.Range("I132").Value = 661.2368
$ws.Range("K132").Value = 1983.7104
$ws.Range("M132").Value = 546.2896000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 534
$ws.Range("I6").Value = 534
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 534
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -361
$ws.Range("N6").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 749514.8
$ws.Range("I32").Value = 786422.4399999999
$ws.Range("K32").Value = 786422.4399999999
$ws.Range("M32").Value = -786135.4399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1875.1052
$ws.Range("I110").Value = 1703.375
$ws.Range("K110").Value = 1703.375
$ws.Range("M110").Value = 341.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 44477.4
$ws.Range("J58").Value = 51669.5
$ws.Range("L58").Value = 51669.5
$ws.Range("N58").Value = -52257.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 96150
$ws.Range("J59").Value = 96150
$ws.Range("L59").Value = 96150
$ws.Range("N59").Value = -97844

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 9351.736999999999
$ws.Range("I105").Value = 8052.846
$ws.Range("J105").Value = 12166
$ws.Range("K105").Value = 8052.846
$ws.Range("L105").Value = 12166
$ws.Range("M105").Value = -6305.846
$ws.Range("N105").Value = -15660

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4186.222
$ws.Range("I16").Value = 4186.222
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 4186.222
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3899.222
$ws.Range("N16").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2289864.2
$ws.Range("I31").Value = 2507812.8
$ws.Range("J31").Value = 1403.5
$ws.Range("K31").Value = 2507812.8
$ws.Range("L31").Value = 1403.5
$ws.Range("M31").Value = -2507517.8
$ws.Range("N31").Value = -1993.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2289864.2
$ws.Range("I34").Value = 2507812.8
$ws.Range("J34").Value = 1403.5
$ws.Range("K34").Value = 2507812.8
$ws.Range("L34").Value = 1403.5
$ws.Range("M34").Value = -2507610.8
$ws.Range("N34").Value = -1807.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 37587740
$ws.Range("I58").Value = 83337340
$ws.Range("K58").Value = 83337340
$ws.Range("M58").Value = -83337137

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 4186.222
$ws.Range("I113").Value = 4186.222
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4186.222
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2016.222
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 6532.0415
$ws.Range("I132").Value = 6809.9546
$ws.Range("K132").Value = 20429.8638
$ws.Range("M132").Value = -17899.8638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 41429.5
$ws.Range("J133").Value = 40633.715
$ws.Range("L133").Value = 40633.715
$ws.Range("N133").Value = -45693.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2324.1155
$ws.Range("I134").Value = 1975.0869
$ws.Range("K134").Value = 5925.2607
$ws.Range("M134").Value = -3390.2607

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 37587740
$ws.Range("I136").Value = 83337340
$ws.Range("K136").Value = 250012020
$ws.Range("M136").Value = -250009470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 33.714287
$ws.Range("I40").Value = 31.2
$ws.Range("J40").Value = 40
$ws.Range("K40").Value = 124.8
$ws.Range("L40").Value = 160
$ws.Range("M40").Value = -55.8
$ws.Range("N40").Value = -298

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1077084.4
$ws.Range("I122").Value = 3227144.2
$ws.Range("K122").Value = 29044297.8
$ws.Range("M122").Value = -29041847.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").Value = 46326
$ws.Range("J125").Value = 46326
$ws.Range("L125").Value = 46326
$ws.Range("N125").Value = -51246

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3262.0386
$ws.Range("I22").Value = 2765.889
$ws.Range("J22").Value = 3524.7058
$ws.Range("K22").Value = 2765.889
$ws.Range("L22").Value = 3524.7058
$ws.Range("M22").Value = -2470.889
$ws.Range("N22").Value = -4114.7058

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3262.0386
$ws.Range("I27").Value = 2765.889
$ws.Range("J27").Value = 3524.7058
$ws.Range("K27").Value = 2765.889
$ws.Range("L27").Value = 3524.7058
$ws.Range("M27").Value = -2658.889
$ws.Range("N27").Value = -3738.7058

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5245
$ws.Range("I40").Value = 4660.1665
$ws.Range("J40").Value = 6999.5
$ws.Range("K40").Value = 4660.1665
$ws.Range("L40").Value = 6999.5
$ws.Range("M40").Value = -4524.1665
$ws.Range("N40").Value = -7271.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 267333.34
$ws.Range("J43").Value = 267333.34
$ws.Range("L43").Value = 267333.34
$ws.Range("N43").Value = -267719.34

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3895.44
$ws.Range("I46").Value = 1071.4286
$ws.Range("K46").Value = 1071.4286
$ws.Range("M46").Value = -883.4286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1762.3478
$ws.Range("I82").Value = 1479.3529
$ws.Range("K82").Value = 1479.3529
$ws.Range("M82").Value = -1118.3529

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1762.3478
$ws.Range("I85").Value = 1479.3529
$ws.Range("K85").Value = 1479.3529
$ws.Range("M85").Value = -231.3529000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3623.8572
$ws.Range("I93").Value = 2870.1875
$ws.Range("K93").Value = 2870.1875
$ws.Range("M93").Value = -1622.1875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 44999.668
$ws.Range("J124").Value = 44999.668
$ws.Range("L124").Value = 44999.668
$ws.Range("N124").Value = -54819.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2086101
$ws.Range("I132").Value = 4169042.5
$ws.Range("J132").Value = 3159.5625
$ws.Range("K132").Value = 12507127.5
$ws.Range("L132").Value = 9478.6875
$ws.Range("M132").Value = -12504597.5
$ws.Range("N132").Value = -14538.6875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 82247.75
$ws.Range("J140").Value = 86664
$ws.Range("L140").Value = 86664
$ws.Range("N140").Value = -97024

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5954734
$ws.Range("I132").Value = 7577871
$ws.Range("K132").Value = 22733613
$ws.Range("M132").Value = -22731083
